$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge the three runs around "к примеру," in the "Массив" paragraph into
#    a single run and drop the proofErr gramStart/gramEnd markers.
#    (The exact character span was pre-computed from the original document:
#    it is the run sequence "Массив ... Другие структуры данных, " +
#    "к примеру," + " стеки и очереди, являются производными от ".)
# ---------------------------------------------------------------------------
$mergeFind = $d.Content.Find
$mergeFind.ClearFormatting()
$mergeFind.Text = "Массив — одна из самых простых"
$mergeFind.Execute() | Out-Null
$mergeStart = $mergeFind.Parent.Start

$endFind = $d.Content.Find
$endFind.ClearFormatting()
$endFind.Text = "являются производными от "
$endFind.Execute() | Out-Null
$mergeEnd = $endFind.Parent.End

$mergeRange = $d.Range($mergeStart, $mergeEnd)
$newText = "Массив — одна из самых простых и часто применяемых структур данных. Другие структуры данных, к примеру, стеки и очереди, являются производными от "

# Force an actual content change (a same-text assignment is a no-op in this
# engine and would leave the old run/proofErr split untouched), then set the
# final text in a second pass.
$mergeRange.Text = "#"
$mergeRange2 = $d.Range($mergeStart, $mergeStart + 1)
$mergeRange2.Text = $newText

# ---------------------------------------------------------------------------
# 2) At the very end of the document: keep the last paragraph ("... именно
#    словарь нам надо реализовать.") untouched, then append a new paragraph
#    holding a page break, then a final empty paragraph. Wrap the whole
#    document body (from its start) in a bookmark named "_Hlk184832737"
#    that ends right after the page-break paragraph.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$tailRange = $lastPara.Range
$tailRange.Collapse(0)
$tailRange.InsertParagraphAfter()

$breakParaRange = $d.Paragraphs.Last.Range
$breakParaRange.Collapse(1)
$breakParaRange.InsertBreak(7)

$breakPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$breakPara.Format.SpaceAfter = 8
$breakPara.Format.LineSpacingRule = 5
$breakPara.Format.LineSpacing = 12.95
$breakPara.Alignment = 0

$finalPara = $d.Paragraphs.Last
$finalRange = $finalPara.Range
$finalRange.InsertAfter("#")
$cleanupRange = $d.Range($finalRange.Start, $finalRange.Start + 1)
$cleanupRange.Delete()

$bookmarkEndPos = $breakPara.Range.End
$bookmarkRange = $d.Range(0, $bookmarkEndPos)
$d.Bookmarks.Add("_Hlk184832737", $bookmarkRange)
